{"js": "const replacements = [\n  [\"37\u00d740=1480\", \"66\u00d777=5082\"],\n  [\"87\u00d713=1131\", \"64\u00d778=4992\"],\n  [\"91\u00d793=8463\", \"19\u00d748=912\"],\n  [\"27\u00d758=1566\", \"57\u00d743=2451\"],\n  [\"17\u00d779=1343\", \"21\u00d793=1953\"],\n  [\"89\u00d775=6675\", \"71\u00d776=5396\"],\n  [\"94\u00d718=1692\", \"47\u00d714=658\"],\n  [\"97\u00d753=5141\", \"44\u00d769=3036\"],\n  [\"44\u00d761=2684\", \"36\u00d793=3348\"],\n  [\"53\u00d749=2597\", \"25\u00d760=1500\"],\n  [\"64\u00d768=4352\", \"69\u00d771=4899\"],\n  [\"31\u00d792=2852\", \"69\u00d740=2760\"],\n  [\"76\u00d776=5776\", \"18\u00d780=1440\"],\n  [\"85\u00d736=3060\", \"29\u00d761=1769\"],\n  [\"80\u00d731=2480\", \"73\u00d788=6424\"],\n  [\"29\u00d773=2117\", \"99\u00d732=3168\"],\n  [\"21\u00d792=1932\", \"11\u00d798=1078\"],\n  [\"55\u00d795=5225\", \"56\u00d747=2632\"],\n  [\"63\u00d793=5859\", \"30\u00d753=1590\"],\n  [\"81\u00d754=4374\", \"20\u00d743=860\"],\n  [\"11\u00d724=264\", \"58\u00d774=4292\"],\n  [\"48\u00d780=3840\", \"39\u00d779=3081\"],\n  [\"76\u00d785=6460\", \"33\u00d716=528\"],\n  [\"20\u00d713=260\", \"86\u00d714=1204\"],\n  [\"65\u00d713=845\", \"67\u00d725=1675\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"37\u00d740=1480\"\n$find.Replacement.Text = \"66\u00d777=5082\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"87\u00d713=1131\"\n$find.Replacement.Text = \"64\u00d778=4992\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"91\u00d793=8463\"\n$find.Replacement.Text = \"19\u00d748=912\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"27\u00d758=1566\"\n$find.Replacement.Text = \"57\u00d743=2451\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"17\u00d779=1343\"\n$find.Replacement.Text = \"21\u00d793=1953\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"89\u00d775=6675\"\n$find.Replacement.Text = \"71\u00d776=5396\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"94\u00d718=1692\"\n$find.Replacement.Text = \"47\u00d714=658\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"97\u00d753=5141\"\n$find.Replacement.Text = \"44\u00d769=3036\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"44\u00d761=2684\"\n$find.Replacement.Text = \"36\u00d793=3348\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"53\u00d749=2597\"\n$find.Replacement.Text = \"25\u00d760=1500\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"64\u00d768=4352\"\n$find.Replacement.Text = \"69\u00d771=4899\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"31\u00d792=2852\"\n$find.Replacement.Text = \"69\u00d740=2760\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"76\u00d776=5776\"\n$find.Replacement.Text = \"18\u00d780=1440\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"85\u00d736=3060\"\n$find.Replacement.Text = \"29\u00d761=1769\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"80\u00d731=2480\"\n$find.Replacement.Text = \"73\u00d788=6424\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"29\u00d773=2117\"\n$find.Replacement.Text = \"99\u00d732=3168\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"21\u00d792=1932\"\n$find.Replacement.Text = \"11\u00d798=1078\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"55\u00d795=5225\"\n$find.Replacement.Text = \"56\u00d747=2632\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"63\u00d793=5859\"\n$find.Replacement.Text = \"30\u00d753=1590\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"81\u00d754=4374\"\n$find.Replacement.Text = \"20\u00d743=860\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"11\u00d724=264\"\n$find.Replacement.Text = \"58\u00d774=4292\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"48\u00d780=3840\"\n$find.Replacement.Text = \"39\u00d779=3081\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"76\u00d785=6460\"\n$find.Replacement.Text = \"33\u00d716=528\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"20\u00d713=260\"\n$find.Replacement.Text = \"86\u00d714=1204\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"65\u00d713=845\"\n$find.Replacement.Text = \"67\u00d725=1675\"\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $false, $null, 2) | Out-Null\n"}
